# Re-label the per-category detail rows so each line item states its
# category (e.g. "New nominations" -> "Civilian, New nominations"), and
# turn the old generic "Summary" block into two explicit total rows,
# shifting what used to be row 40 out of existence (table now ends at 39).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Civilian ---
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

# --- Other Civilian ---
$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Returned to White House "

# --- Air Force ---
$ws.Range("A17").Value = "     Air Force, New nominations"
$ws.Range("A18").Value = "     Air Force, Confirmed "
$ws.Range("A19").Value = "     Air Force, Withdrawn "
$ws.Range("A20").Value = "     Air Force, Returned to White House "

# --- Army ---
$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Confirmed "
$ws.Range("A24").Value = "     Army, Returned to White House "

# --- Navy ---
$ws.Range("A26").Value = "     Navy, New nominations"
$ws.Range("A27").Value = "     Navy, Confirmed "
$ws.Range("A28").Value = "     Navy, Withdrawn "
$ws.Range("A29").Value = "     Navy, Returned to White House "

# --- Marine Corps ---
$ws.Range("A31").Value = "     Marine Corps, New nominations"
$ws.Range("A32").Value = "     Marine Corps, Confirmed "
$ws.Range("A33").Value = "     Marine Corps, Returned to White House "

# --- Summary block: replace the bare "Summary" header row with a real
# "Total new nominations" figure (19074, matching what used to be on the
# row below), then relabel/shift the remaining total rows up by one.
$ws.Range("A34").Value = "Total new nominations"
$ws.Range("B34").Value = 19074
# B34 is a brand-new number cell; give it the same "#,##0" look as the
# other totals (style 3) by copying formats from an existing cell that
# already carries that style.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B34").PasteSpecial(-4122) | Out-Null

$ws.Range("A35").Value = "Total carryover nominations"
# B35 keeps its value (0) and style - unchanged

$ws.Range("A36").Value = "Total confirmed "
$ws.Range("B36").Value = 17328
# B36 keeps its style (3) - unchanged

$ws.Range("A37").Value = "Total unconfirmed "
$ws.Range("B37").Value = 2
# B37 switches from the thousands-separator style (3) to the plain
# numeric style (2) used elsewhere for small counts.
$ws.Range("B38").Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null

$ws.Range("A38").Value = "Total withdrawn "
$ws.Range("B38").Value = 11
# B38 keeps its style (2) - unchanged

$ws.Range("A39").Value = "Total returned to the White House "
$ws.Range("B39").Value = 1733
# B39 switches from the plain numeric style (2) to the thousands-
# separator style (3), matching the old row 40 it replaces.
$ws.Range("B36").Copy() | Out-Null
$ws.Range("B39").PasteSpecial(-4122) | Out-Null

# The old row 40 ("Total returned to the White House " / 1733) has now
# been folded into row 39 above, so drop the now-duplicate trailing row.
$ws.Rows(40).Delete()
